$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: split the run that holds the "-m "Message here in quotes" ... "
# sentence of the "Committing" bullet into several runs, adding the
# w:proofErr grammar markers that Word's grammar checker inserted around
# the quoted dash and the doubled "always  use" phrase.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Committing*use*git commit*My message*committing new changes.*") {
        $rng = $p.Range
        $xmlFragment = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="1E53A379" w14:textId="534DD117" w:rsidR="00D64431" w:rsidRDefault="00D64431" w:rsidP="007713CE">
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="2"/>
<w:numId w:val="1"/>
</w:numPr>
</w:pPr>
<w:r><w:t xml:space="preserve">Committing – use “git commit </w:t></w:r>
<w:r><w:t xml:space="preserve">-m “Message here in quotes” </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>“ –</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>always  use</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> git commit -m “My message” ----- while you are committing new changes.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
        [void]$rng.InsertXML($xmlFragment)
        break
    }
}

# ---------------------------------------------------------------------------
# Part 2: append two new bullet paragraphs (same list, ilvl=1) after the
# "When making a commit ..." paragraph, describing "Git Log" and
# "Git add .".
# ---------------------------------------------------------------------------
$insertPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$xmlFragment2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Git </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>Log :</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> It retrieves the information for the logs of the commits for a repository</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Git </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>add .</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> ---- it stages all the files with the changes for the commit. If there are multiple changes to the files tracked by git then it will add all the files for the commit. </w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
[void]$insertPoint.InsertXML($xmlFragment2)

Write-Host "Edit complete"
